# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 130
$sheet1.Range("F3").Value = 324
$sheet1.Range("F5").Value = 1606
$sheet1.Range("F7").Value = 2153
$sheet1.Range("F9").Value = 282
$sheet1.Range("F10").Value = 106
$sheet1.Range("F11").Value = 4811
$sheet1.Range("F17").Value = 168
$sheet1.Range("F20").Value = 112
$sheet1.Range("F21").Value = 3744
$sheet1.Range("F22").Value = 690
$sheet1.Range("F23").Value = 618
$sheet1.Range("F27").Value = 114
$sheet1.Range("F33").Value = 20
$sheet1.Range("F34").Value = 863
$sheet1.Range("F35").Value = 2360
$sheet1.Range("F36").Value = 423

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 130
$sheet4.Range("F3").Value = 324
$sheet4.Range("F5").Value = 1606
$sheet4.Range("F7").Value = 2153
$sheet4.Range("F9").Value = 282
$sheet4.Range("F10").Value = 106
$sheet4.Range("F11").Value = 4811
$sheet4.Range("F17").Value = 168
$sheet4.Range("F20").Value = 112
$sheet4.Range("F21").Value = 3744
$sheet4.Range("F22").Value = 690
$sheet4.Range("F23").Value = 618
$sheet4.Range("F27").Value = 114
$sheet4.Range("F34").Value = 20
$sheet4.Range("F35").Value = 863
$sheet4.Range("F36").Value = 2360
$sheet4.Range("F37").Value = 423
